# NDF_Generator: environment verifications & optional installations
# - Hide the "Options" sheet (kept around for data-validation lists, but no
#   longer meant to be browsed directly by the user).
# - Replace the placeholder "test" line item on the dashboard with the real
#   basketball-club purchase list (ballons, panier, chaussures, maillots,
#   shorts) including their quantities / unit prices / totals.

$wb = $excel.ActiveWorkbook

$dashboard = $wb.Worksheets.Item("Tableau de bord")
$options   = $wb.Worksheets.Item("Options")

# --- Fill in the real purchase list on the dashboard (columns A-D) --------
# Quantité | Référence | Prix unitaire | Prix total
$dashboard.Range("A2").Value = 12
$dashboard.Range("B2").Value = "Ballons de basket"
$dashboard.Range("C2").Value = 20
$dashboard.Range("D2").Value = 240

$dashboard.Range("A3").Value = 1
$dashboard.Range("B3").Value = "Panier de basket"
$dashboard.Range("C3").Value = 325
$dashboard.Range("D3").Value = 325

$dashboard.Range("A4").Value = 1
$dashboard.Range("B4").Value = "Paire de chaussure"
$dashboard.Range("C4").Value = 150
$dashboard.Range("D4").Value = 150

$dashboard.Range("A5").Value = 25
$dashboard.Range("B5").Value = "Maillots"
$dashboard.Range("C5").Value = 10
$dashboard.Range("D5").Value = 250

$dashboard.Range("A6").Value = 25
$dashboard.Range("B6").Value = "Short"
$dashboard.Range("C6").Value = 15
$dashboard.Range("D6").Value = 375

# --- Hide the "Options" sheet ---------------------------------------------
$options.Visible = $false
